$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended after the existing last row (340).
# Columns: A=Empresa, B=Categoria, C=Produto, D=Mes, E=porcem_repasse
$rows = @(
    @("Imóveis", "Serviço", "Avaliação", 1, 0.12),
    @("Imóveis", "Serviço", "Avaliação", 2, 0),
    @("Imóveis", "Serviço", "Lançamento", 1, 0.12),
    @("Imóveis", "Serviço", "Lançamento", 2, 0),
    @("Imóveis", "Serviço", "Consultoria e Incorporação", 1, 0.1),
    @("Imóveis", "Serviço", "Consultoria e Incorporação", 2, 0),
    @("Imóveis", "Imóveis Prontos", "Indicando Comprador ou Vendedor", 1, 0.12),
    @("Imóveis", "Imóveis Prontos", "Indicando Comprador ou Vendedor", 2, 0),
    @("Imóveis", "Imóveis Prontos", "Indicando Comprador e Vendedor", 1, 0.24),
    @("Imóveis", "Imóveis Prontos", "Indicando Comprador e Vendedor", 2, 0)
)

$startRow = 341
$rowIndex = $startRow

# Written in the same per-row order the source file was authored in
# (Categoria, then Produto, then Empresa, then Mes/porcem_repasse) so new
# shared-string entries land at the same indices as the original edit.
foreach ($row in $rows) {
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex = $rowIndex + 1
}

# Column D on the new rows reuses the existing bordered style already used
# by column D throughout the sheet (copy formats only, from the row above
# the new block, so no duplicate style entries are created).
$ws.Range("D340").Copy()
$ws.Range("D341:D350").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Trailing formatting-only row: an empty, underlined cell in column E.
$lastRow = $rowIndex
$ws.Cells.Item($lastRow, 5).Font.Underline = 1

# Reflect the navigation state left behind at save time (best effort).
$excel.ActiveWindow.ScrollRow = 315
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E336").Select()
